# Add the "2021-2022" worksheet (copy of the "2020-2021" tax table) and
# make it the active sheet, matching the author's commit:
#   "Implemented get_year_prompt function - Updated tax rates for year 2021-2022"

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2020-2021")

# New sheet goes after the existing one, in tab order.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "2021-2022"

# Copy the tax-bracket table (headers + 5 brackets) across, preserving the
# shared-string reuse ("Tax_Bracket", "Lower_Limit", ..., "INF").
$src.Range("A1:E6").Copy()
$newSheet.Range("A1").PasteSpecial()

# Match the source sheet's column widths.
$newSheet.Columns.Item(1).ColumnWidth = 10.666666666666666
$newSheet.Columns.Item(2).ColumnWidth = 9
$newSheet.Columns.Item(3).ColumnWidth = 21.333333333333336
$newSheet.Columns.Item(4).ColumnWidth = 11.666666666666666
$newSheet.Columns.Item(5).ColumnWidth = 16

# Restore the current selection on the new sheet.
$newSheet.Range("C11").Select() | Out-Null

$excel.CutCopyMode = $false
